$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Status column (D) values and apply the matching builtin
# "Good"/"Bad"/"Neutral" cell styles (the Excel Cell Styles gallery),
# processed top-to-bottom so the new fonts/fills/styles are registered
# in the same order Excel itself would create them (Good, Bad, Neutral).

$ws.Range("D7").Value2 = "Done"
$ws.Range("D7").Style = "Good"

$ws.Range("D8").Value2 = "Done"
$ws.Range("D8").Style = "Good"

$ws.Range("D9").Value2 = "Done"
$ws.Range("D9").Style = "Good"

$ws.Range("D10").Value2 = "Done"
$ws.Range("D10").Style = "Good"

$ws.Range("D11").Value2 = "Not Implemented"
$ws.Range("D11").Style = "Bad"

$ws.Range("D12").Value2 = "Done"
$ws.Range("D12").Style = "Good"

$ws.Range("D13").Value2 = "Done"
$ws.Range("D13").Style = "Good"

$ws.Range("D14").Value2 = "Not Implemented"
$ws.Range("D14").Style = "Bad"

$ws.Range("D15").Value2 = "In Progress"
$ws.Range("D15").Style = "Neutral"

$ws.Range("D16").Value2 = "Done"
$ws.Range("D16").Style = "Good"

$ws.Range("D17").Value2 = "In Progress"
$ws.Range("D17").Style = "Neutral"

$ws.Range("D18").Value2 = "Not Implemented"
$ws.Range("D18").Style = "Bad"

$ws.Range("D19").Value2 = "Not Implemented"
$ws.Range("D19").Style = "Bad"

$ws.Range("D20").Value2 = "Not Implemented"
$ws.Range("D20").Style = "Bad"

$ws.Range("D21").Value2 = "Not Implemented"
$ws.Range("D21").Style = "Bad"

$ws.Range("D22").Value2 = "Not Implemented"
$ws.Range("D22").Style = "Bad"

$ws.Range("D23").Value2 = "Not Implemented"
$ws.Range("D23").Style = "Bad"

# Update the saved cursor/selection position for the sheet view.
$ws.Range("F18").Select() | Out-Null

Write-Host "Applied status styling and selection updates"
